# Auto-update draw results: append the 2025-10-16 Pick 4 draw as a new
# row at the bottom of the results table (row 30), matching the nightly
# "Auto-update draw results" job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 30

# Force text storage for every new cell (the whole A:E table is kept as
# text - see the sheet's numberStoredAsText ignoredError - so dates like
# "2025-10-16" and numeric-looking phases like "251016" are not silently
# reinterpreted as a date serial / number by the COM layer).
$ws.Range("A" + $row + ":E" + $row).NumberFormat = "@"

$ws.Range("A" + $row).Value = "2025-10-16"
$ws.Range("B" + $row).Value = "Pick 4"
$ws.Range("C" + $row).Value = "251016"
$ws.Range("D" + $row).Value = "5-6-1-1"
$ws.Range("E" + $row).Value = "2025-10-16T21:38:22.616+04:00"
